{"js": "// Fix the stray closing bracket in the merge-field placeholder that was\n// typed as \"{navnAvsender]\" (mismatched bracket) so it reads correctly\n// as \"{navnAvsender}\" in the \"Postmottaket til ...\" paragraph.\nconst body = context.document.body;\n\nconst results = body.search(\"{navnAvsender]\", { matchCase: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{navnAvsender}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix the stray closing bracket in the merge-field placeholder that was\n# typed as \"{navnAvsender]\" (mismatched bracket) so it reads correctly\n# as \"{navnAvsender}\" in the \"Postmottaket til ...\" paragraph.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"{navnAvsender]\"\n$find.Replacement.Text = \"{navnAvsender}\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n"}
